$d = $word.ActiveDocument

$replacements = @(
    @("91×19=", "50×63="),
    @("96×73=", "87×53="),
    @("38×28=", "70×43="),
    @("24×13=", "91×79="),
    @("72×47=", "78×25="),
    @("70×45=", "80×11="),
    @("34×80=", "16×45="),
    @("87×84=", "98×81="),
    @("96×69=", "43×53="),
    @("20×58=", "22×19="),
    @("76×84=", "24×20="),
    @("45×63=", "47×53="),
    @("96×58=", "85×86="),
    @("64×85=", "25×29="),
    @("34×63=", "69×95="),
    @("57×54=", "59×33="),
    @("84×22=", "93×33="),
    @("37×22=", "87×94="),
    @("85×93=", "35×98="),
    @("80×69=", "49×27="),
    @("91×72=", "70×99="),
    @("72×50=", "52×28="),
    @("45×28=", "67×36="),
    @("34×73=", "81×33="),
    @("86×73=", "58×33=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
